# Test_FridgesSorting.xlsx edit
# - Rename the two header labels on the sheet:
#     "Category"    -> "category"
#     "How to sort" -> "howToSort"
# - Move the active cell / selection from A4 to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "category"
$ws.Range("B1").Value = "howToSort"

$ws.Range("B2").Select()
